$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 449.104309
$ws.Range("H2").Value = 1347.312927
$ws.Range("I2").Value = 0.9710020245482639
$ws.Range("J2").Value = 0.9710020245482639
$ws.Range("M2").Value = 60.813934
$ws.Range("N2").Value = 182.441802
$ws.Range("O2").Value = 0.6840634102070431
$ws.Range("P2").Value = 0.6840634102070431
$ws.Range("Q2").Value = 27311.7998066416
$ws.Range("R2").Value = 245806.1982597744
$ws.Range("S2").Value = 0.6642269562304284
$ws.Range("T2").Value = 0.6642269562304284
$ws.Range("G3").Value = 449.104309
$ws.Range("H3").Value = 1347.312927
$ws.Range("I3").Value = 0.9710020245482639
$ws.Range("J3").Value = 0.9710020245482639
$ws.Range("M3").Value = 2.823821666666666
$ws.Range("N3").Value = 8.471464999999998
$ws.Range("O3").Value = 0.0317636592810545
$ws.Range("P3").Value = 0.0317636592810545
$ws.Range("Q3").Value = 1268.190478347561
$ws.Range("R3").Value = 11413.71430512805
$ws.Range("S3").Value = 0.03084257746896517
$ws.Range("T3").Value = 0.03084257746896517
$ws.Range("G4").Value = 449.104309
$ws.Range("H4").Value = 1347.312927
$ws.Range("I4").Value = 0.9710020245482639
$ws.Range("J4").Value = 0.9710020245482639
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.640208
$ws.Range("N4").Value = 1.920624
$ws.Range("O4").Value = 0.007201357302782462
$ws.Range("P4").Value = 0.007201357302782463
$ws.Range("Q4").Value = 287.520171456272
$ws.Range("R4").Value = 2587.681543106448
$ws.Range("S4").Value = 0.006992532520497195
$ws.Range("T4").Value = 0.006992532520497196
$ws.Range("G5").Value = 449.104309
$ws.Range("H5").Value = 1347.312927
$ws.Range("I5").Value = 0.9710020245482639
$ws.Range("J5").Value = 0.9710020245482639
$ws.Range("M5").Value = 24.623055
$ws.Range("N5").Value = 73.86916500000001
$ws.Range("O5").Value = 0.2769715732091199
$ws.Range("P5").Value = 0.2769715732091199
$ws.Range("Q5").Value = 11058.320101244
$ws.Range("R5").Value = 99524.88091119597
$ws.Range("S5").Value = 0.2689399583283731
$ws.Range("T5").Value = 0.2689399583283731
$ws.Range("I6").Value = 0.01131353526791385
$ws.Range("J6").Value = 0.01131353526791385
$ws.Range("M6").Value = 60.813934
$ws.Range("N6").Value = 182.441802
$ws.Range("O6").Value = 0.6840634102070431
$ws.Range("P6").Value = 0.6840634102070431
$ws.Range("Q6").Value = 318.22076837213
$ws.Range("R6").Value = 2863.98691534917
$ws.Range("S6").Value = 0.007739175516866798
$ws.Range("T6").Value = 0.007739175516866799
$ws.Range("I7").Value = 0.01131353526791385
$ws.Range("J7").Value = 0.01131353526791385
$ws.Range("M7").Value = 2.823821666666666
$ws.Range("N7").Value = 8.471464999999998
$ws.Range("O7").Value = 0.0317636592810545
$ws.Range("P7").Value = 0.0317636592810545
$ws.Range("S7").Value = 0.000359359279514209
$ws.Range("T7").Value = 0.000359359279514209
$ws.Range("I8").Value = 0.01131353526791385
$ws.Range("J8").Value = 0.01131353526791385
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.640208
$ws.Range("N8").Value = 1.920624
$ws.Range("O8").Value = 0.007201357302782462
$ws.Range("P8").Value = 0.007201357302782463
$ws.Range("Q8").Value = 3.35001320056
$ws.Range("R8").Value = 30.15011880504
$ws.Range("S8").Value = [double]"8.147280982187831E-05"
$ws.Range("T8").Value = [double]"8.147280982187832E-05"
$ws.Range("I9").Value = 0.01131353526791385
$ws.Range("J9").Value = 0.01131353526791385
$ws.Range("M9").Value = 24.623055
$ws.Range("N9").Value = 73.86916500000001
$ws.Range("O9").Value = 0.2769715732091199
$ws.Range("P9").Value = 0.2769715732091199
$ws.Range("Q9").Value = 128.844936783225
$ws.Range("R9").Value = 1159.604431049025
$ws.Range("S9").Value = 0.00313352766171096
$ws.Range("T9").Value = 0.00313352766171096
$ws.Range("G10").Value = 5.266527
$ws.Range("H10").Value = 15.799581
$ws.Range("I10").Value = 0.0113866829528418
$ws.Range("J10").Value = 0.0113866829528418
$ws.Range("M10").Value = 60.813934
$ws.Range("N10").Value = 182.441802
$ws.Range("O10").Value = 0.6840634102070431
$ws.Range("P10").Value = 0.6840634102070431
$ws.Range("Q10").Value = 320.278225387218
$ws.Range("R10").Value = 2882.504028484962
$ws.Range("S10").Value = 0.007789213171667363
$ws.Range("T10").Value = 0.007789213171667363
$ws.Range("G11").Value = 5.266527
$ws.Range("H11").Value = 15.799581
$ws.Range("I11").Value = 0.0113866829528418
$ws.Range("J11").Value = 0.0113866829528418
$ws.Range("M11").Value = 2.823821666666666
$ws.Range("N11").Value = 8.471464999999998
$ws.Range("O11").Value = 0.0317636592810545
$ws.Range("P11").Value = 0.0317636592810545
$ws.Range("Q11").Value = 14.871733050685
$ws.Range("R11").Value = 133.845597456165
$ws.Range("S11").Value = 0.0003616827176554584
$ws.Range("T11").Value = 0.0003616827176554584
$ws.Range("G12").Value = 5.266527
$ws.Range("H12").Value = 15.799581
$ws.Range("I12").Value = 0.0113866829528418
$ws.Range("J12").Value = 0.0113866829528418
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.640208
$ws.Range("N12").Value = 1.920624
$ws.Range("O12").Value = 0.007201357302782462
$ws.Range("P12").Value = 0.007201357302782463
$ws.Range("Q12").Value = 3.371672717616
$ws.Range("R12").Value = 30.345054458544
$ws.Range("S12").Value = [double]"8.199957243691585E-05"
$ws.Range("T12").Value = [double]"8.199957243691585E-05"
$ws.Range("G13").Value = 5.266527
$ws.Range("H13").Value = 15.799581
$ws.Range("I13").Value = 0.0113866829528418
$ws.Range("J13").Value = 0.0113866829528418
$ws.Range("M13").Value = 24.623055
$ws.Range("N13").Value = 73.86916500000001
$ws.Range("O13").Value = 0.2769715732091199
$ws.Range("P13").Value = 0.2769715732091199
$ws.Range("Q13").Value = 129.677983979985
$ws.Range("R13").Value = 1167.101855819865
$ws.Range("S13").Value = 0.003153787491082059
$ws.Range("T13").Value = 0.003153787491082059
$ws.Range("G14").Value = 2.912815666666667
$ws.Range("H14").Value = 8.738447000000001
$ws.Range("I14").Value = 0.006297757230980464
$ws.Range("J14").Value = 0.006297757230980464
$ws.Range("M14").Value = 60.813934
$ws.Range("N14").Value = 182.441802
$ws.Range("O14").Value = 0.6840634102070431
$ws.Range("P14").Value = 0.6840634102070431
$ws.Range("Q14").Value = 177.1397797068327
$ws.Range("R14").Value = 1594.258017361494
$ws.Range("S14").Value = 0.004308065288080561
$ws.Range("T14").Value = 0.004308065288080561
$ws.Range("G15").Value = 2.912815666666667
$ws.Range("H15").Value = 8.738447000000001
$ws.Range("I15").Value = 0.006297757230980464
$ws.Range("J15").Value = 0.006297757230980464
$ws.Range("M15").Value = 2.823821666666666
$ws.Range("N15").Value = 8.471464999999998
$ws.Range("O15").Value = 0.0317636592810545
$ws.Range("P15").Value = 0.0317636592810545
$ws.Range("Q15").Value = 8.225271990539444
$ws.Range("R15").Value = 74.02744791485499
$ws.Range("S15").Value = 0.0002000398149196607
$ws.Range("T15").Value = 0.0002000398149196607
$ws.Range("G16").Value = 2.912815666666667
$ws.Range("H16").Value = 8.738447000000001
$ws.Range("I16").Value = 0.006297757230980464
$ws.Range("J16").Value = 0.006297757230980464
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.640208
$ws.Range("N16").Value = 1.920624
$ws.Range("O16").Value = 0.007201357302782462
$ws.Range("P16").Value = 0.007201357302782463
$ws.Range("Q16").Value = 1.864807892325334
$ws.Range("R16").Value = 16.783271030928
$ws.Range("S16").Value = [double]"4.535240002647222E-05"
$ws.Range("T16").Value = [double]"4.535240002647222E-05"
$ws.Range("G17").Value = 2.912815666666667
$ws.Range("H17").Value = 8.738447000000001
$ws.Range("I17").Value = 0.006297757230980464
$ws.Range("J17").Value = 0.006297757230980464
$ws.Range("M17").Value = 24.623055
$ws.Range("N17").Value = 73.86916500000001
$ws.Range("O17").Value = 0.2769715732091199
$ws.Range("P17").Value = 0.2769715732091199
$ws.Range("Q17").Value = 71.72242036519502
$ws.Range("R17").Value = 645.5017832867551
$ws.Range("S17").Value = 0.00174429972795377
$ws.Range("T17").Value = 0.00174429972795377
